$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 ("Rule" column of the last rule row) text changes from "R40" to "1".
# Prefix with an apostrophe so Excel stores it as text (shared string "1")
# instead of silently converting it to the number 1.
$ws.Range("B11").Value = "'1"
